$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.85"
$ws.Range("E2").Value = "'7.98%"
$ws.Range("D3").Value = "'31.85"
$ws.Range("E3").Value = "'8.96%"
$ws.Range("D4").Value = "'5.353"
$ws.Range("D5").Value = "'0.07571"
$ws.Range("E5").Value = "'13.19%"
$ws.Range("D6").Value = "'7.844"
$ws.Range("E6").Value = "'6.81%"
$ws.Range("D7").Value = "'3.716"
$ws.Range("E7").Value = "'9.04%"
$ws.Range("D8").Value = "'1.625"
$ws.Range("E8").Value = "'18.39%"
$ws.Range("D9").Value = "'0.9208"
$ws.Range("E9").Value = "'0.15%"
$ws.Range("D10").Value = "'0.01682"
$ws.Range("E10").Value = "'2,502.85%"
$ws.Range("D11").Value = "'0.1715"
$ws.Range("E11").Value = "'7.71%"
$ws.Range("D12").Value = "'0.07625"
$ws.Range("E12").Value = "'11.90%"
$ws.Range("D13").Value = "'0.08205"
$ws.Range("E13").Value = "'7.86%"
$ws.Range("D14").Value = "'0.03032"
$ws.Range("E14").Value = "'3.43%"
$ws.Range("D15").Value = "'0.09890"
$ws.Range("E15").Value = "'10.16%"
$ws.Range("D16").Value = "'0.001539"
$ws.Range("E16").Value = "'-2.21%"
$ws.Range("D17").Value = "'0.04559"
$ws.Range("E17").Value = "'1.20%"
$ws.Range("D18").Value = "'0.006586"
$ws.Range("E18").Value = "'4.44%"
$ws.Range("D19").Value = "'3.494"
$ws.Range("E19").Value = "'1.27%"
$ws.Range("E20").Value = "'1.09%"
$ws.Range("D21").Value = "'0.3312"
$ws.Range("E21").Value = "'3.08%"
$ws.Range("D22").Value = "'0.1337"
$ws.Range("E22").Value = "'2.12%"
$ws.Range("D23").Value = "'4.216"
$ws.Range("E23").Value = "'3.61%"
$ws.Range("D24").Value = "'0.1629"
$ws.Range("E24").Value = "'2.94%"
$ws.Range("D25").Value = "'0.001229"
$ws.Range("E25").Value = "'3.31%"
$ws.Range("D26").Value = "'0.004488"
$ws.Range("E26").Value = "'9.18%"
$ws.Range("D27").Value = "'0.0001301"
$ws.Range("E27").Value = "'8.45%"
$ws.Range("D28").Value = "'0.0001742"
$ws.Range("E28").Value = "'7.64%"
$ws.Range("D40").Value = "'0.04577"
$ws.Range("E40").Value = "'7.28%"
$ws.Range("D41").Value = "'0.007207"
$ws.Range("E42").Value = "'10.42%"
$ws.Range("D43").Value = "'0.002261"
$ws.Range("E43").Value = "'1.46%"
$ws.Range("D44").Value = "'0.01405"
$ws.Range("E44").Value = "'4.81%"
$ws.Range("D45").Value = "'0.00006158"
$ws.Range("E45").Value = "'7.98%"
$ws.Range("D46").Value = "'1.893"
$ws.Range("E46").Value = "'-3.83%"
$ws.Range("D47").Value = "'0.01301"
$ws.Range("E47").Value = "'-0.50%"
